$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link/Price/Volume cells store plain text (prices like "36.592.84" use dots as
# thousands separators, not decimal points, so force Text format before writing so
# Excel does not reinterpret them as numbers/dates).
$cellValues = [ordered]@{
    "D2" = "36.592.84"
    "E2" = "  -0.94%  "
    "D3" = "2.058.32"
    "E3" = "  +0.97%  "
    "D4" = "1.00"
    "E4" = "  -0.01%  "
    "D5" = "243.27"
    "E5" = "  -0.75%  "
    "D6" = "0.668"
    "E6" = "  +1.43%  "
    "E7" = "  +0.08%  "
    "D8" = "54.52"
    "E8" = "  -6.97%  "
    "D9" = "58.63"
    "E9" = "  -0.86%  "
    "D10" = "0.362"
    "E10" = "  -3.96%  "
    "E11" = "  -2.08%  "
    "E12" = "  -2.99%  "
    "D13" = "0.934"
    "E13" = "  +5.71%  "
    "D14" = "14.74"
    "E14" = "  -4.31%  "
    "D15" = "2.361.04"
    "E15" = "  +1.20%  "
    "D16" = "5.43"
    "E16" = "  -3.76%  "
    "D17" = "2.062.00"
    "E17" = "  +1.48%  "
    "D18" = "36.503.58"
    "E18" = "  -1.07%  "
    "D19" = "16.76"
    "E19" = "  -7.78%  "
    "D20" = "71.99"
    "E20" = "  -2.17%  "
    "E21" = "  -3.09%  "
    "D22" = "238.09"
    "E22" = "  +1.00%  "
    "E23" = "  -2.16%  "
    "E24" = "  -0.12%  "
    "E25" = "  -3.89%  "
    "E26" = "  +1.55%  "
    "D27" = "9.30"
    "E27" = "  -3.44%  "
    "D28" = "164.81"
    "E28" = "  -2.58%  "
    "D29" = "20.11"
    "E29" = "  +0.93%  "
    "E30" = "  -1.21%  "
    "E31" = "  +8.85%  "
    "D32" = "5.05"
    "E32" = "  -7.87%  "
    "D33" = "4.49"
    "E33" = "  -5.04%  "
    "D34" = "0.0598"
    "E34" = "  -2.37%  "
    "E35" = "  +0.18%  "
    "E36" = "  -0.63%  "
    "E37" = "  -2.08%  "
    "D38" = "0.0820"
    "E38" = "  -5.18%  "
    "E39" = "  -4.07%  "
    "E40" = "  -5.48%  "
    "E41" = "  -2.21%  "
    "E42" = "  -7.55%  "
    "E43" = "  -2.50%  "
    "D44" = "94.33"
    "E44" = "  -3.00%  "
    "B45" = "Maker"
    "C45" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D45" = "1.422.26"
    "E45" = "  +9.92%  "
    "B46" = "Cronos"
    "C46" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D46" = "0.0911"
    "E46" = "  -5.75%  "
    "B47" = "InjectiveProtocol"
    "C47" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D47" = "16.00"
    "E47" = "  -5.54%  "
    "B48" = "FraxShare"
    "C48" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D48" = "7.55"
    "E48" = "  +12.87%  "
    "D49" = "2.86"
    "E49" = "  +0.65%  "
    "E50" = "  -2.44%  "
    "D51" = "2.248.21"
    "E51" = "  +1.28%  "
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
}
